$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-shuffle the Gantt rows:
#    - "Programming" moves from row 6 to row 7
#    - "Testing" (row 9) becomes "Testing/Bug-Fixing" with new dates
#    - "Running" (row 10) becomes "Running/Extras" with new dates
#    - "Extras" (row 11) is removed (row goes blank)
#    - Row 8 ("Writing Code") keeps its label but gets a new end date
# ---------------------------------------------------------------------------

# Row 6: was "Programming" 43499-43562 -> now blank (keep blank styled cells)
$ws.Range("A6:D6").ClearContents()

# Row 7: was blank -> now "Programming" 43499-43562
$ws.Range("A7").Value = "Programming"
$ws.Range("B7").Value = 43499
$ws.Range("C7").Value = 43562
$ws.Range("D7").Formula = "=-(B7-C7)"
$ws.Range("D7").Style = "Normal"

# Row 8: "Writing Code" end date 43538 -> 43525 (duration 39 -> 26, auto via formula)
$ws.Range("C8").Value = 43525

# Row 10: "Running" -> "Running/Extras", 43538-43539 -> 43538-43562
$ws.Range("A10").Value = "Running/Extras"
$ws.Range("C10").Value = 43562

# Row 9: "Testing" -> "Testing/Bug-Fixing", 43538-43562 -> 43525-43538
$ws.Range("A9").Value = "Testing/Bug-Fixing"
$ws.Range("B9").Value = 43525
$ws.Range("C9").Value = 43538

# Row 11: was "Extras" 43538-43539 -> now blank
$ws.Range("A11:D11").ClearContents()

# ---------------------------------------------------------------------------
# 2. Chart cosmetic tweaks that are controllable from the COM surface
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# Narrower gap between bars
$cg = $chart.ChartGroups(1)
$cg.GapWidth = 41

# Reposition / resize the chart to its new anchor (from col B ~ to col N)
$co.Left = 210.3934
$co.Top = 10.6382
$co.Width = 693.6509
$co.Height = 491.7446

# ---------------------------------------------------------------------------
# 3. Selection moves from R10 to B16
# ---------------------------------------------------------------------------
$null = $ws.Range("B16").Select()

Write-Host "Gantt chart updated"
